$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.036036036036036
$ws.Cells.Item(2, 3).Value = 0.000750750750750751
$ws.Cells.Item(2, 4).Value = 0.0015015015015015
$ws.Cells.Item(2, 5).Value = 0.003003003003003
$ws.Cells.Item(2, 6).Value = 0.00225225225225225
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0.990990990990991
$ws.Cells.Item(2, 9).Value = 0.021021021021021
$ws.Cells.Item(2, 10).Value = 0.990990990990991
$ws.Cells.Item(2, 11).Value = 0.0195195195195195
$ws.Cells.Item(2, 12).Value = 0.00600600600600601
$ws.Cells.Item(2, 13).Value = 0.021021021021021
$ws.Cells.Item(2, 14).Value = 0.000750750750750751
$ws.Cells.Item(2, 15).Value = 0.0315315315315315
$ws.Cells.Item(2, 16).Value = 0.0015015015015015
$ws.Cells.Item(2, 17).Value = 0.99024024024024
$ws.Cells.Item(2, 18).Value = 0.00525525525525526
$ws.Cells.Item(2, 19).Value = 0.989489489489489
$ws.Cells.Item(2, 20).Value = 0.048048048048048
$ws.Cells.Item(2, 21).Value = 0.990990990990991
$ws.Cells.Item(2, 22).Value = 0.989489489489489
$ws.Cells.Item(2, 23).Value = 0.00225225225225225
$ws.Cells.Item(2, 24).Value = 0.00225225225225225
$ws.Cells.Item(3, 2).Value = 0.003003003003003
$ws.Cells.Item(3, 3).Value = 0.99024024024024
$ws.Cells.Item(3, 4).Value = 0.992492492492492
$ws.Cells.Item(3, 5).Value = 0.994744744744745
$ws.Cells.Item(3, 6).Value = 0.00225225225225225
$ws.Cells.Item(3, 7).Value = 0.996996996996997
$ws.Cells.Item(3, 8).Value = 0.000750750750750751
$ws.Cells.Item(3, 9).Value = 0.00225225225225225
$ws.Cells.Item(3, 10).Value = 0.00375375375375375
$ws.Cells.Item(3, 11).Value = 0.00225225225225225
$ws.Cells.Item(3, 12).Value = 0.018018018018018
$ws.Cells.Item(3, 13).Value = 0.972972972972973
$ws.Cells.Item(3, 14).Value = 0.00225225225225225
$ws.Cells.Item(3, 15).Value = 0.00375375375375375
$ws.Cells.Item(3, 16).Value = 0.996996996996997
$ws.Cells.Item(3, 17).Value = 0.00600600600600601
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 0.003003003003003
$ws.Cells.Item(3, 20).Value = 0.0015015015015015
$ws.Cells.Item(3, 21).Value = 0.00225225225225225
$ws.Cells.Item(3, 22).Value = 0.00600600600600601
$ws.Cells.Item(3, 23).Value = 0.00525525525525526
$ws.Cells.Item(3, 24).Value = 0.984234234234234
$ws.Cells.Item(4, 2).Value = 0.952702702702703
$ws.Cells.Item(4, 3).Value = 0.00375375375375375
$ws.Cells.Item(4, 4).Value = 0.000750750750750751
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0.000750750750750751
$ws.Cells.Item(4, 7).Value = 0.000750750750750751
$ws.Cells.Item(4, 8).Value = 0.0045045045045045
$ws.Cells.Item(4, 9).Value = 0.972222222222222
$ws.Cells.Item(4, 10).Value = 0.00375375375375375
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0.972222222222222
$ws.Cells.Item(4, 13).Value = 0.003003003003003
$ws.Cells.Item(4, 14).Value = 0.0045045045045045
$ws.Cells.Item(4, 15).Value = 0.963963963963964
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 0.003003003003003
$ws.Cells.Item(4, 19).Value = 0.00525525525525526
$ws.Cells.Item(4, 20).Value = 0.00225225225225225
$ws.Cells.Item(4, 21).Value = 0.0015015015015015
$ws.Cells.Item(4, 22).Value = 0.003003003003003
$ws.Cells.Item(4, 23).Value = 0.987237237237237
$ws.Cells.Item(4, 24).Value = 0.00825825825825826
$ws.Cells.Item(5, 2).Value = 0.00600600600600601
$ws.Cells.Item(5, 3).Value = 0.00525525525525526
$ws.Cells.Item(5, 4).Value = 0.00525525525525526
$ws.Cells.Item(5, 5).Value = 0.00225225225225225
$ws.Cells.Item(5, 6).Value = 0.994744744744745
$ws.Cells.Item(5, 7).Value = 0.00225225225225225
$ws.Cells.Item(5, 8).Value = 0.00375375375375375
$ws.Cells.Item(5, 9).Value = 0.0045045045045045
$ws.Cells.Item(5, 10).Value = 0.0015015015015015
$ws.Cells.Item(5, 11).Value = 0.978228228228228
$ws.Cells.Item(5, 12).Value = 0.00375375375375375
$ws.Cells.Item(5, 13).Value = 0.003003003003003
$ws.Cells.Item(5, 14).Value = 0.992492492492492
$ws.Cells.Item(5, 15).Value = 0.000750750750750751
$ws.Cells.Item(5, 16).Value = 0.0015015015015015
$ws.Cells.Item(5, 17).Value = 0.00375375375375375
$ws.Cells.Item(5, 18).Value = 0.991741741741742
$ws.Cells.Item(5, 19).Value = 0.00225225225225225
$ws.Cells.Item(5, 20).Value = 0.948198198198198
$ws.Cells.Item(5, 21).Value = 0.00525525525525526
$ws.Cells.Item(5, 22).Value = 0.0015015015015015
$ws.Cells.Item(5, 23).Value = 0.00525525525525526
$ws.Cells.Item(5, 24).Value = 0.00525525525525526
